$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 275
$ws.Range("B13").Value = -550
$ws.Range("B17").Value = -0.551
$ws.Range("B18").Value = 0.835

$ws.Range("B18").Select()
